# maj template comment a la fin
# Move the "Comment" column (J) to the end of the block (Y), shifting the
# columns in between (K:Y) one position to the left, for every header/meta
# row of the template (rows 1-5: field name, French label, type hint,
# format note, example).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startCol = 10   # J
$endCol   = 25   # Y

for ($row = 1; $row -le 5; $row++) {
    $values = @()
    for ($col = $startCol; $col -le $endCol; $col++) {
        $values += $ws.Cells.Item($row, $col).Value2
    }

    # rotate left by one: first element (Comment column) goes to the end
    $first = $values[0]
    for ($i = 0; $i -lt ($values.Length - 1); $i++) {
        $values[$i] = $values[$i + 1]
    }
    $values[$values.Length - 1] = $first

    for ($i = 0; $i -lt $values.Length; $i++) {
        $ws.Cells.Item($row, $startCol + $i).Value = $values[$i]
    }
}
